$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells
$ws.Range("C126").Value = 16
$ws.Range("C1083").Value = 7
$ws.Range("C1280").Value = 4
$ws.Range("C1297").Value = 15
$ws.Range("C1301").Value = 13
$ws.Range("C1302").Value = 17
$ws.Range("C1305").Value = 15
$ws.Range("C1306").Value = 4
$ws.Range("C1307").Value = 5
$ws.Range("C1308").Value = 7
$ws.Range("C1309").Value = 12

$ws.Range("B1310").Value = "20-29"
$ws.Range("C1310").Value = 1

$ws.Range("B1311").Value = "50-59"
$ws.Range("C1311").Value = 2

# New rows 1312-1324
$newRows = @(
    @(44226, "60-69", 2),
    @(44226, "70-79", 13),
    @(44226, "80+", 9),
    @(44227, "40-49", 3),
    @(44227, "50-59", 2),
    @(44227, "60-69", 7),
    @(44227, "70-79", 7),
    @(44227, "80+", 6),
    @(44228, "40-49", 1),
    @(44228, "50-59", 1),
    @(44228, "60-69", 1),
    @(44228, "70-79", 2),
    @(44228, "80+", 6)
)

$r = 1312
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
